$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the teacher name in columns B/C (no label in A) is
# removed; everything below shifts up by one row.
$ws.Rows(13).Delete()

# A handful of cells keep their (now shifted) row/column position but get
# their contents swapped for different values.
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "5840535 - Messias Borges Silva"
$ws.Range("C18").Value = "5840535 - Messias Borges Silva"

$recuperacao = "NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação"
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao
